# trafo_id -> gridnode_id refactor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J header: rename "trafo_id" to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Move/restore the active cell selection to E6
$ws.Range("E6").Select()
